$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the results table: the "train_data_control" rows (previously
# rows 8-9) move to the top (rows 2-3); the remaining blocks shift down
# in their original relative order.
$data = @(
    @("train_data_control", "k-NearestNeighbour_classifier", 0.81113727095915655, 0.79614872458487684, 0.81113727095915655, 0.78879480622315157),
    @("train_data_control", "DecisionTree_classifier",       0.76998669259903774, 0.77376309992929382, 0.76998669259903774, 0.77177381193877159),
    @("kNN_imputed_10",     "k-NearestNeighbour_classifier", 0.80663322755655642, 0.79114315683906911, 0.80663322755655642, 0.79250085765419731),
    @("kNN_imputed_10",     "DecisionTree_classifier",       0.70907974204115054, 0.74608020309512013, 0.70907974204115054, 0.72264749886069268),
    @("kNN_imputed_40",     "k-NearestNeighbour_classifier", 0.56658818712253045, 0.75494946409264829, 0.56658818712253045, 0.59484576263802624),
    @("kNN_imputed_40",     "DecisionTree_classifier",       0.54611526256525744, 0.72284536196831894, 0.54611526256525744, 0.57673703679745625),
    @("kNN_imputed_70",     "k-NearestNeighbour_classifier", 0.26164397584194898, 0.72185401505440516, 0.26164397584194898, 0.14568881689230639),
    @("kNN_imputed_70",     "DecisionTree_classifier",       0.39082812979834169, 0.68432286173780987, 0.39082812979834169, 0.39160898070338468)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# View changes: zoomed-in Normal view, with the selection moved from J13 to J5.
$excel.ActiveWindow.Zoom = 130
[void]$ws.Range("J5").Select()
